$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("H2").Value = 0.004002571105957031

# Row 3
$ws.Range("C3").Value = 644.41
$ws.Range("D3").Value = 469.7
$ws.Range("E3").Value = 174.71
$ws.Range("F3").Value = 1114.11
$ws.Range("G3").Value = 557.05
$ws.Range("H3").Value = 1.652802228927612

# Row 4
$ws.Range("H4").Value = 0.003083229064941406

# Row 5
$ws.Range("C5").Value = 515.71
$ws.Range("D5").Value = 514
$ws.Range("E5").Value = 1.71
$ws.Range("F5").Value = 1029.71
$ws.Range("G5").Value = 514.86
$ws.Range("H5").Value = 1.614155054092407

# Row 6
$ws.Range("H6").Value = 0.006043195724487305

# Row 7
$ws.Range("C7").Value = 555.53
$ws.Range("D7").Value = 551.71
$ws.Range("E7").Value = 3.82
$ws.Range("F7").Value = 1661.68
$ws.Range("G7").Value = 553.89
$ws.Range("H7").Value = 3.694030523300171

# Row 8
$ws.Range("H8").Value = 0.0110476016998291

# Row 9
$ws.Range("C9").Value = 693.99
$ws.Range("D9").Value = 448.81
$ws.Range("E9").Value = 245.18
$ws.Range("F9").Value = 2334.93
$ws.Range("G9").Value = 583.73
$ws.Range("H9").Value = 6.131538867950439

# Row 10
$ws.Range("H10").Value = 0.01109099388122559

# Row 11
$ws.Range("C11").Value = 540.41
$ws.Range("D11").Value = 529.09
$ws.Range("E11").Value = 11.32
$ws.Range("F11").Value = 2145.14
$ws.Range("G11").Value = 536.28
$ws.Range("H11").Value = 6.246362447738647
